$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly price sheet gained a new (most-recent) reporting date.
# Insert two fresh rows right after the header block of this product
# (before the current first data pair at row 57) so every existing
# weekly pair shifts down by one week (2 rows: "Primera"/"Segunda").
$ws.Range("A57:A58").EntireRow.Insert()

# Populate the two newly inserted rows with the new week's data.
# Row 57 - Calidad "Primera"
$ws.Cells.Item(57,1).Value = 1
$ws.Cells.Item(57,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(57,3).Value = "Arica y Parinacota"
$ws.Cells.Item(57,4).Value = 44413
$ws.Cells.Item(57,5).Value = 15
$ws.Cells.Item(57,6).Value = 100114014
$ws.Cells.Item(57,7).Value = "Betarraga"
$ws.Cells.Item(57,8).Value = "Sin especificar"
$ws.Cells.Item(57,9).Value = "Primera"
$ws.Cells.Item(57,10).Value = 700
$ws.Cells.Item(57,11).Value = 500
$ws.Cells.Item(57,12).Value = 600
$ws.Cells.Item(57,13).Value = 550
$ws.Cells.Item(57,14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(57,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(57,16).Value = 138
$ws.Cells.Item(57,17).Value = 4
$ws.Cells.Item(57,18).Value = "Hortaliza"

# Row 58 - Calidad "Segunda"
$ws.Cells.Item(58,1).Value = 1
$ws.Cells.Item(58,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(58,3).Value = "Arica y Parinacota"
$ws.Cells.Item(58,4).Value = 44413
$ws.Cells.Item(58,5).Value = 15
$ws.Cells.Item(58,6).Value = 100114014
$ws.Cells.Item(58,7).Value = "Betarraga"
$ws.Cells.Item(58,8).Value = "Sin especificar"
$ws.Cells.Item(58,9).Value = "Segunda"
$ws.Cells.Item(58,10).Value = 900
$ws.Cells.Item(58,11).Value = 500
$ws.Cells.Item(58,12).Value = 600
$ws.Cells.Item(58,13).Value = 550
$ws.Cells.Item(58,14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(58,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(58,16).Value = 110
$ws.Cells.Item(58,17).Value = 5
$ws.Cells.Item(58,18).Value = "Hortaliza"
